$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.195.76'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '3.573.79'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '199.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '588.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.624'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.214'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.640'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.05'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000299'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.55'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '693.34'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +15.71%  '
$ws.Range('D15').Value = '4.156.46'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').Value = '70.321.37'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '12.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.82%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.593.88'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.66%  '
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.988'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '110.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.48'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.98'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.38%  '
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.91'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.33'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.41'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.03'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.16'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.113'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '63.42'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = '0.0₃0840'
$ws.Range('E36').Value = '  +2.17%  '
$ws.Range('D37').Value = '3.799.57'
$ws.Range('E37').Value = '  -1.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.33%  '
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '509.34'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.14%  '
$ws.Range('E41').Value = '  -7.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.20'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.380'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.57%  '
$ws.Range('E44').Value = '  +2.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0467'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.140'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.59'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('E50').Value = '  -0.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.80'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +20.82%  '
